# "update week 3 thurs" -- rework the tentative-schedule table (weeks 3-9)
# and resize the columns on the grade table and the schedule table.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Column width tweaks (w:tblGrid / w:gridCol widths, expressed in
#    points on the COM Width property == twentieths-of-a-point / 20).
# ---------------------------------------------------------------------

$gradeTable = $d.Tables.Item(1)
$gradeTable.Columns.Item(1).Width = 99.0    # 1980 twips
$gradeTable.Columns.Item(2).Width = 159.5   # 3190 twips

$scheduleTable = $d.Tables.Item(2)
$scheduleTable.Columns.Item(1).Width = 82.5   # 1650 twips
$scheduleTable.Columns.Item(2).Width = 148.5  # 2970 twips
$scheduleTable.Columns.Item(3).Width = 82.5   # 1650 twips
$scheduleTable.Columns.Item(4).Width = 82.5   # 1650 twips

# ---------------------------------------------------------------------
# 2. Re-write the week 3 - week 9 rows of the schedule table. Each
#    week's content effectively shifts: week 3 gets a brand new topic,
#    and weeks 4-9 pick up the prior week's (adjusted) content.
# ---------------------------------------------------------------------

$scheduleTable = $d.Tables.Item(2)

# Week 3 (row 4)
$scheduleTable.Cell(4, 2).Range.Text = "Descriptive statistics and graphical summaries"
$scheduleTable.Cell(4, 3).Range.Text = "1.4 – 1.6"

# Week 4 (row 5)
$scheduleTable.Cell(5, 2).Range.Text = "Foundations for inference"
$scheduleTable.Cell(5, 3).Range.Text = "4.1 – 4.2"
$scheduleTable.Cell(5, 4).Range.Text = "HW2"

# Week 5 (row 6)
$scheduleTable.Cell(6, 2).Range.Text = "One- and two-sample inference for numerical data"
$scheduleTable.Cell(6, 3).Range.Text = "4.3; 5.1 — 5.3"
$scheduleTable.Cell(6, 4).Range.Text = "Test 1"

# Week 6 (row 7)
$scheduleTable.Cell(7, 2).Range.Text = "Nonparametric inference"
$scheduleTable.Cell(7, 3).Range.Text = "TBD"
$scheduleTable.Cell(7, 4).Range.Text = "HW3"

# Week 7 (row 8)
$scheduleTable.Cell(8, 2).Range.Text = "Comparing many means with analysis of variance"
$scheduleTable.Cell(8, 3).Range.Text = "5.5"
$scheduleTable.Cell(8, 4).Range.Text = "HW4"

# Week 8 (row 9)
$scheduleTable.Cell(9, 2).Range.Text = "Inference for categorical data"
$scheduleTable.Cell(9, 3).Range.Text = "8.1 — 8.4"
$scheduleTable.Cell(9, 4).Range.Text = "Test 2"

# Week 9 (row 10) -- only the assignment column changes (HW6 -> HW5)
$scheduleTable.Cell(10, 4).Range.Text = "HW5"

# ---------------------------------------------------------------------
# 3. Strike through the "HW2" assignment for week 3 (row 4) -- the
#    homework text stays the same but is now marked as struck-out,
#    without touching the paragraph mark's own run formatting.
# ---------------------------------------------------------------------

$hw2Cell = $scheduleTable.Cell(4, 4)
$hw2CellRange = $hw2Cell.Range
$hw2Range = $d.Range($hw2CellRange.Start, $hw2CellRange.End - 1)
$hw2Range.Font.StrikeThrough = 1
